# "Remove guids from m_logs detail"
#
# The survey sheet had a duplicate/guid-bearing "begin screen" / "end screen"
# pair (rows 6 and 9) wrapping the date/notes fields (rows 10-11) before the
# real "end screen" (row 12). Collapse that by deleting the redundant rows
# 6-9, which shifts the date/notes/end-screen rows up to 6-8.
#
# It also updates the "settings" sheet so the form's instance_name setting
# points at the "date" field instead of the removed "date_serviced" field.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$survey.Rows("6:9").Delete()

# Reset the view: selection lands on the new first empty-ish row beneath the
# data, and the sheet scrolls back to showing column A.
$survey.Range("B11").Select() | Out-Null

$settings = $wb.Worksheets.Item("settings")
$settings.Range("B6").Value = "date"
